# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# to match the target diff for 'Sheets/Ultros_Profits.xlsx' (scheduled profit-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 764.6957
$ws.Range("I28").Value = 626.35297
$ws.Range("K28").Value = 626.35297
$ws.Range("M28").Value = -141.35297
$ws.Range("H43").Value = 3960.2
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 3950.25
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 3950.25
$ws.Range("M43").Value = -3931
$ws.Range("N43").Value = -4088.25
$ws.Range("H116").Value = 10807.25
$ws.Range("J116").Value = 11117
$ws.Range("L116").Value = 11117
$ws.Range("N116").Value = -18001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2127.6956
$ws.Range("I74").Value = 2250.7058
$ws.Range("J74").Value = 1779.1666
$ws.Range("K74").Value = 2250.7058
$ws.Range("L74").Value = 1779.1666
$ws.Range("M74").Value = -1376.7058
$ws.Range("N74").Value = -3527.1666
$ws.Range("H77").Value = 2127.6956
$ws.Range("I77").Value = 2250.7058
$ws.Range("J77").Value = 1779.1666
$ws.Range("K77").Value = 11253.529
$ws.Range("L77").Value = 8895.833000000001
$ws.Range("M77").Value = -6885.529
$ws.Range("N77").Value = -17631.833
$ws.Range("H132").Value = 2332.5881
$ws.Range("I132").Value = 2332.5881
$ws.Range("K132").Value = 6997.7643
$ws.Range("M132").Value = -4467.7643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 55558764
$ws.Range("I86").Value = 83335650
$ws.Range("K86").Value = 83335650
$ws.Range("M86").Value = -83334527
$ws.Range("H89").Value = 55558764
$ws.Range("I89").Value = 83335650
$ws.Range("K89").Value = 416678250
$ws.Range("M89").Value = -416672634
$ws.Range("H94").Value = 2281.5
$ws.Range("I94").Value = 1947.52
$ws.Range("K94").Value = 1947.52
$ws.Range("M94").Value = -1496.52

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1448.919
$ws.Range("I31").Value = 1290.4
$ws.Range("J31").Value = 2128.2856
$ws.Range("K31").Value = 1290.4
$ws.Range("L31").Value = 2128.2856
$ws.Range("M31").Value = -995.4000000000001
$ws.Range("N31").Value = -2718.2856
$ws.Range("H34").Value = 1448.919
$ws.Range("I34").Value = 1290.4
$ws.Range("J34").Value = 2128.2856
$ws.Range("K34").Value = 1290.4
$ws.Range("L34").Value = 2128.2856
$ws.Range("M34").Value = -1088.4
$ws.Range("N34").Value = -2532.2856
$ws.Range("H107").Value = 7937603.5
$ws.Range("I107").Value = 11905779
$ws.Range("J107").Value = 1252.8334
$ws.Range("K107").Value = 11905779
$ws.Range("L107").Value = 1252.8334
$ws.Range("M107").Value = -11903859
$ws.Range("N107").Value = -5092.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1061.5862
$ws.Range("I5").Value = 999.7406999999999
$ws.Range("K5").Value = 2999.2221
$ws.Range("M5").Value = -2887.2221
$ws.Range("H22").Value = 380.77777
$ws.Range("I22").Value = 403.375
$ws.Range("K22").Value = 1210.125
$ws.Range("M22").Value = -1041.125
$ws.Range("H27").Value = 380.77777
$ws.Range("I27").Value = 403.375
$ws.Range("K27").Value = 1210.125
$ws.Range("M27").Value = -1108.125
$ws.Range("H68").Value = 2390.7273
$ws.Range("I68").Value = 959.6
$ws.Range("K68").Value = 2878.8
$ws.Range("M68").Value = -2067.8
$ws.Range("H69").Value = 2416.6667
$ws.Range("I69").Value = 1958.3334
$ws.Range("K69").Value = 5875.0002
$ws.Range("M69").Value = -5064.0002
$ws.Range("H71").Value = 2390.7273
$ws.Range("I71").Value = 959.6
$ws.Range("K71").Value = 8636.4
$ws.Range("M71").Value = -4580.4
$ws.Range("H72").Value = 2416.6667
$ws.Range("I72").Value = 1958.3334
$ws.Range("K72").Value = 17625.0006
$ws.Range("M72").Value = -13569.0006
$ws.Range("H94").Value = 4000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H125").Value = 5500
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 2248
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2248
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 6744
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -16744
$ws.Range("H130").Value = 3500000
$ws.Range("I130").Value = 3500000
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 10500000
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -10494980
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 3267.0625
$ws.Range("I131").Value = 1743.3
$ws.Range("J131").Value = 5806.6665
$ws.Range("K131").Value = 5229.9
$ws.Range("L131").Value = 17419.9995
$ws.Range("M131").Value = -189.8999999999996
$ws.Range("N131").Value = -27499.9995
$ws.Range("H135").Value = 1061.5862
$ws.Range("I135").Value = 999.7406999999999
$ws.Range("K135").Value = 8997.666299999999
$ws.Range("M135").Value = -6462.666299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10115.5
$ws.Range("J33").Value = 10115.5
$ws.Range("L33").Value = 10115.5
$ws.Range("N33").Value = -10619.5
$ws.Range("H55").Value = 30
$ws.Range("I55").Value = 30
$ws.Range("K55").Value = 30
$ws.Range("M55").Value = 297
$ws.Range("H132").Value = 7801.5674
$ws.Range("I132").Value = 7166.4116
$ws.Range("K132").Value = 21499.2348
$ws.Range("M132").Value = -18969.2348

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 90911250
$ws.Range("I82").Value = 125002160
$ws.Range("J82").Value = 2143.6667
$ws.Range("K82").Value = 125002160
$ws.Range("L82").Value = 2143.6667
$ws.Range("M82").Value = -125001799
$ws.Range("N82").Value = -2865.6667
$ws.Range("H85").Value = 90911250
$ws.Range("I85").Value = 125002160
$ws.Range("J85").Value = 2143.6667
$ws.Range("K85").Value = 125002160
$ws.Range("L85").Value = 2143.6667
$ws.Range("M85").Value = -125000912
$ws.Range("N85").Value = -4639.6667
$ws.Range("H122").Value = 4763.4814
$ws.Range("I122").Value = 4240.15
$ws.Range("J122").Value = 6258.7144
$ws.Range("K122").Value = 12720.45
$ws.Range("L122").Value = 18776.1432
$ws.Range("M122").Value = -10270.45
$ws.Range("N122").Value = -23676.1432
$ws.Range("H132").Value = 2028.4482
$ws.Range("I132").Value = 1743.75
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5231.25
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -2701.25
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3425
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3566.6667
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3566.6667
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4814.6667
$ws.Range("H65").Value = 3425
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3566.6667
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 17833.3335
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -24073.3335
$ws.Range("H132").Value = 3950
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -170
$ws.Range("N132").Value = -26060

